# Update logBook worksheet with entries through 28th June (rows 30 and 31),
# matching the authoritative diff:
#  - two new shared strings
#  - two new rows (30, 31) with same formatting as row 29
#  - E29's shared formula group extended logically (D-C) down through E31
#  - selection moved to G32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 29 down into the two new rows so columns keep
# the same number formats / alignment (date, time, wrap text, etc.)
$ws.Range("A29:G29").Copy()
$ws.Range("A30:G30").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A31:G31").PasteSpecial(-4122)   # xlPasteFormats

# Row 30 -> Sno 29, 28-Jun-2022, 05:00-05:30, Code, "1. deeplabv3 paper review"
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 44740
$ws.Cells.Item(30, 3).Value = 0.20833333333333334
$ws.Cells.Item(30, 4).Value = 0.22916666666666666
$ws.Cells.Item(30, 6).Value = "Code"
$ws.Cells.Item(30, 7).Value = "1. deeplabv3 paper review"

# Row 31 -> Sno 30, 28-Jun-2022, 20:00-20:30, Code, "1. deeplab_v3_starter nb"
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 44740
$ws.Cells.Item(31, 3).Value = 0.83333333333333337
$ws.Cells.Item(31, 4).Value = 0.85416666666666663
$ws.Cells.Item(31, 6).Value = "Code"
$ws.Cells.Item(31, 7).Value = "1. deeplab_v3_starter nb"

# Extend the Time column formula (Time = endTime - startTime) into the new rows
$ws.Range("E30:E31").Formula = "=D30-C30"

# Move the active selection the way the workbook was left after the edit
$ws.Range("G32").Select()
